$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Swap the Fecha (date) values between row 2 and row 4
$ws.Range("D2").Value = 44557
$ws.Range("D4").Value = 44547

# Swap the Volumen values between row 2 and row 4
$ws.Range("J2").Value = 400
$ws.Range("J4").Value = 200
